# Update "Azerbaijan Premier League" results sheet:
#  - row 136 (existing match, id=134) gets corrected/updated odds data
#    (including 3 brand-new columns H/I/J that didn't exist on that row before)
#  - 3 new matches are appended as rows 137-139 (id=135,136,137)
# Data below reflects the full, final state of each cell in rows 136-139,
# columns A..AC (in order). $null marks a column that must stay empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

$row136 = @(134, 7011628, "Azerbaijan Premier League", "Azerbaijan Premier League", 45380.39583333334, "FK Gabala", "Neftchi Baku", 0, 1, "A", 3.5, 4, 1.727, 2.9, 3.75, 1.95, 0.5, 1.75, 1.95, 2.5, 1.95, 1.85, -1, -1, 0.95, -1, 0.95, -1, 0.8500000000000001)
$row137 = @(135, 7011631, "Azerbaijan Premier League", "Azerbaijan Premier League", 45380.5, "FK Kapaz", "FK Qarabag", 1, 6, "A", 6, 5.5, 1.3, 7, 7, 1.25, 2, 1.85, 1.95, 3.5, 1.85, 1.95, -1, -1, 0.25, -1, 0.95, 0.8500000000000001, -1)
$row138 = @(136, 7011630, "Azerbaijan Premier League", "Azerbaijan Premier League", 45381.39583333334, "Sabail FC", "PFK Turan Tovuz", 1, 1, "D", 2.3, 3.2, 2.75, 2.5, 3.25, 2.5, 0, 1.9, 1.9, 2.5, 1.975, 1.825, -1, 2.25, -1, 0, -0, -1, 0.825)
$row139 = @(137, 7011629, "Azerbaijan Premier League", "Azerbaijan Premier League", 45381.5, "Sabah", "Zira IK", 0, 1, "A", 2.1, 3.2, 3.1, 2.05, 3, 3.5, -0.25, 1.85, 1.95, 2, 1.9, 1.9, -1, -1, 2.5, -1, 0.95, -1, 0.8999999999999999)

$allRows = @($row136, $row137, $row138, $row139)
$rowNums = @(136, 137, 138, 139)

for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $r = $rowNums[$i]
    $data = $allRows[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $val = $data[$j]
        if ($val -ne $null) {
            $ws.Cells.Item($r, $j + 1).Value = $val
        }
    }
}

# The new rows (137-139) need the same look as the rest of the table:
# column A is bold/bordered/centered ("id" style) and column E is the
# custom date/time number format. Copy *only* the formatting (not the
# values we just wrote) from the template row 136, which already has
# the correct styles, onto the freshly created rows.
$ws.Range("A136").Copy()
$ws.Range("A137:A139").PasteSpecial(-4122)
$ws.Range("E136").Copy()
$ws.Range("E137:E139").PasteSpecial(-4122)

Write-Host "done"